# Weekly fruit/vegetable price update: insert a new daily-price record for
# "Choclo" (Choclero, Primera, Región de Arica y Parinacota) dated
# 2021-11-09, shifting the existing tail of the table down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 635; everything currently at/after row 635
# (through the old last row 661) shifts down to 636..662, and the sheet's
# used range/dimension grows to A1:R662 automatically.
$ws.Rows(635).Insert()

# Populate the newly inserted row 635 with the new record.
$ws.Range('A635').Value = 6
$ws.Range('B635').Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range('C635').Value = 'Metropolitana'
$ws.Range('D635').Value = '2021-11-09'
$ws.Range('E635').Value = 13
$ws.Range('F635').Value = 100112024
$ws.Range('G635').Value = 'Choclo'
$ws.Range('H635').Value = 'Choclero'
$ws.Range('I635').Value = 'Primera'
$ws.Range('J635').Value = 1030
$ws.Range('K635').Value = 28000
$ws.Range('L635').Value = 30000
$ws.Range('M635').Value = 28874
$ws.Range('N635').Value = '$/malla 30 unidades'
$ws.Range('O635').Value = 'Región de Arica y Parinacota'
$ws.Range('P635').Value = 962
$ws.Range('Q635').Value = 30
$ws.Range('R635').Value = 'Hortaliza'
